# Weekly update for "Fruta, Terminal La Palmera de La Serena - Nectarín":
# insert 5 new rows (new price records for Sun Rise / Venus varieties,
# week of 2022-01-24) above the existing data block, which simply shifts
# the prior rows 216-222 down to 221-227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at row 216 (existing rows 216:222 shift down to 221:227)
$ws.Rows("216:220").Insert()

# Columns that are constant across every data row in this block
$mercadoId  = 8
$mercado    = "Terminal La Palmera de La Serena"
$region     = "Coquimbo"
$codreg     = 4
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria  = "Nectarín"
$origen     = "Región de O'Higgins"

# New rows: Date (serial), Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, PrecioKg, KgUnidad
$newRows = @(
    @{ Row=216; Fecha=44585; Variedad="Sun Rise"; Calidad="Primera";  Volumen=20;  PMin=335000; PMax=340000; PProm=337500; Unidad="$/bins (420 kilos)"; PrecioKg=804; KgUnidad=420 },
    @{ Row=217; Fecha=44585; Variedad="Sun Rise"; Calidad="Segunda";  Volumen=16;  PMin=305000; PMax=310000; PProm=307500; Unidad="$/bins (420 kilos)"; PrecioKg=732; KgUnidad=420 },
    @{ Row=218; Fecha=44585; Variedad="Venus";    Calidad="Especial"; Volumen=20;  PMin=395000; PMax=400000; PProm=397500; Unidad="$/bins (420 kilos)"; PrecioKg=946; KgUnidad=420 },
    @{ Row=219; Fecha=44585; Variedad="Venus";    Calidad="Primera";  Volumen=20;  PMin=335000; PMax=340000; PProm=337500; Unidad="$/bins (420 kilos)"; PrecioKg=804; KgUnidad=420 },
    @{ Row=220; Fecha=44585; Variedad="Venus";    Calidad="Segunda";  Volumen=20;  PMin=305000; PMax=310000; PProm=307500; Unidad="$/bins (420 kilos)"; PrecioKg=732; KgUnidad=420 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $mercadoId
    $ws.Range("B$rowNum").Value = $mercado
    $ws.Range("C$rowNum").Value = $region
    $ws.Range("D$rowNum").Value = $r.Fecha
    $ws.Range("E$rowNum").Value = $codreg
    $ws.Range("F$rowNum").Value = $tipo
    $ws.Range("G$rowNum").Value = $productoId
    $ws.Range("H$rowNum").Value = $producto
    $ws.Range("I$rowNum").Value = $categoriaId
    $ws.Range("J$rowNum").Value = $categoria
    $ws.Range("K$rowNum").Value = $r.Variedad
    $ws.Range("L$rowNum").Value = $r.Calidad
    $ws.Range("M$rowNum").Value = $r.Volumen
    $ws.Range("N$rowNum").Value = $r.PMin
    $ws.Range("O$rowNum").Value = $r.PMax
    $ws.Range("P$rowNum").Value = $r.PProm
    $ws.Range("Q$rowNum").Value = $r.Unidad
    $ws.Range("R$rowNum").Value = $origen
    $ws.Range("S$rowNum").Value = $r.PrecioKg
    $ws.Range("T$rowNum").Value = $r.KgUnidad
}
